# Update countries & provincias Spain
#
# This script reproduces the source diff against paises.xlsx:
#  1. Re-ranks a few country rows so that the country name shown in each
#     row changes (the underlying shared-string table gets reordered by
#     Excel automatically once the cell text is (re)assigned):
#       - row 97/98/99   : Honduras / San Marino / Costa de Marfil
#                            -> Costa de Marfil / Honduras / San Marino
#       - row 129/130/131: Monaco / Liechtenstein / Madagascar
#                            -> Madagascar / Monaco / Liechtenstein
#  2. Refreshes the "Datos actualizados..." timestamp cell.
#  3. Refreshes the daily statistics (Casos totales, Nuevos casos,
#     Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes)
#     for the affected rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Country name re-ranking (shared-string shuffle) -------------------

$ws.Range("A97").Value  = "Costa de Marfil"
$ws.Range("A98").Value  = "Honduras"
$ws.Range("A99").Value  = "San Marino"

$ws.Range("A129").Value = "Madagascar"
$ws.Range("A130").Value = "Monaco"
$ws.Range("A131").Value = "Liechtenstein"

# --- 2. Timestamp refresh --------------------------------------------------

$ws.Range("A1").Value = "Datos actualizados a 6 de Abril de 2020 a las 21:52"

# --- 3. Daily statistics refresh -------------------------------------------

# Row 7 (Alemania)
$ws.Range("B7").Value  = 101806
$ws.Range("C7").Value  = 1683
$ws.Range("E7").Value  = 71426
$ws.Range("G7").Value  = 96
$ws.Range("H7").Value  = 1680

# Row 13 (Suiza)
$ws.Range("B13").Value = 21657
$ws.Range("C13").Value = 557
$ws.Range("E13").Value = 12836
$ws.Range("G13").Value = 50
$ws.Range("H13").Value = 765

# Row 34 (Pakistan)
$ws.Range("E34").Value = 3454
$ws.Range("G34").Value = 6
$ws.Range("H34").Value = 53

# Row 51 (Sudafrica)
$ws.Range("B51").Value = 1686
$ws.Range("C51").Value = 31
$ws.Range("E51").Value = 1579
$ws.Range("G51").Value = 1
$ws.Range("H51").Value = 12

# Row 79 (Republica de Macedonia)
$ws.Range("E79").Value = 517
$ws.Range("G79").Value = 5
$ws.Range("H79").Value = 23

# Row 97 (now Costa de Marfil)
$ws.Range("B97").Value = 323
$ws.Range("C97").Value = 62
$ws.Range("D97").Value = 41
$ws.Range("E97").Value = 279
$ws.Range("F97").Value = 0
$ws.Range("H97").Value = 3

# Row 98 (now Honduras)
$ws.Range("B98").Value = 298
$ws.Range("C98").Value = 30
$ws.Range("D98").Value = 6
$ws.Range("E98").Value = 270
$ws.Range("F98").Value = 10
$ws.Range("H98").Value = 22

# Row 99 (now San Marino)
$ws.Range("B99").Value = 266
$ws.Range("D99").Value = 35
$ws.Range("E99").Value = 199
$ws.Range("F99").Value = 14
$ws.Range("H99").Value = 32

# Row 129 (now Madagascar)
$ws.Range("B129").Value = 82
$ws.Range("C129").Value = 10
$ws.Range("D129").Value = 2
$ws.Range("E129").Value = 80
$ws.Range("F129").Value = 6
$ws.Range("H129").Value = 0

# Row 130 (now Monaco)
$ws.Range("C130").Value = 4
$ws.Range("D130").Value = 4
$ws.Range("E130").Value = 72
$ws.Range("F130").Value = 4

# Row 131 (now Liechtenstein)
$ws.Range("B131").Value = 77
$ws.Range("D131").Value = 55
$ws.Range("E131").Value = 21
$ws.Range("F131").Value = 0
$ws.Range("H131").Value = 1

# Row 141 (Mali)
$ws.Range("D141").Value = 9
$ws.Range("E141").Value = 33

# Row 165 (Libia)
$ws.Range("D165").Value = 1
$ws.Range("E165").Value = 16
